# Updates to spreadsheet from Ken.
#
# This script reproduces the edits captured in the commit diff:
#  - fill in missing end_year (C) values of 9999 for several highway rows
#  - fix a typo'd start_year (Union-Turnpike: 1908 -> 1809)
#  - give the Mass-Rails-1875 row (row 40) a "type" of 't' in D39/D40,
#    drop its stray reviewer-question highlight/note, and restore it to
#    plain (unstyled) formatting like its neighbors
#  - fill in the missing start_year for the I93-I95-NH row (row 68), and
#    replace its reviewer-question note with a blank placeholder
#  - remove the reviewer-question note left on the Seaport-Access row
#    (row 112)
#  - move the active selection to D27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Cambridge-Turnpike-1807 -> add end_year ---
$ws.Range("C12").Value = 9999

# --- Row 13: Union-Turnpike -> fix start_year typo, add end_year ---
$ws.Range("B13").Value = 1809
$ws.Range("C13").Value = 9999

# --- Row 14: Route-107-1868 -> add end_year ---
$ws.Range("C14").Value = 9999

# --- Row 16: Route-109-1838 -> add end_year ---
$ws.Range("C16").Value = 9999

# --- Row 17: Route-2-East-1829 -> add end_year ---
$ws.Range("C17").Value = 9999

# --- Row 18: Route-2-West-1830 -> add end_year ---
$ws.Range("C18").Value = 9999

# --- Row 20: Route-1-South-1857 -> add end_year ---
$ws.Range("C20").Value = 9999

# --- Row 22: Route-9-1845 -> add end_year ---
$ws.Range("C22").Value = 9999

# --- Row 23: Newburyport-Turnpike-1805 -> add end_year ---
$ws.Range("C23").Value = 9999

# --- Row 39: Rockport-CR-1861 -> give it a "type" of 't' ---
$ws.Range("D39").Value = "t"

# --- Row 40: Mass-Rails-1875 -> strip the yellow highlight/customFormat,
#     add an end_year, drop the reviewer note ---
$ws.Range("A40").EntireRow.ClearFormats()
$ws.Range("C40").Value = 9999
$ws.Range("E40").Clear()

# --- Row 68: I93-I95-NH -> add start_year, blank out the reviewer note
#     (kept on its existing yellow-highlight style) ---
$ws.Range("B68").Value = 1964
$ws.Range("E68").Value = " "

# --- Row 112: Seaport-Access-2003 -> drop the reviewer note ---
$ws.Range("E112").Clear()

# --- Remove the highlight color from the remaining highlighted rows
#     (68 and 112) so the style reads as a plain/white background
#     instead of yellow ---
$ws.Range("A68:E68").Interior.ThemeColor = 2
$ws.Range("A68:E68").Interior.TintAndShade = 0
$ws.Range("A112:D112").Interior.ThemeColor = 2
$ws.Range("A112:D112").Interior.TintAndShade = 0

# --- Move the active selection ---
$ws.Range("D27").Select()
